$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

# Clear out the old data area first (rows 2-10, columns A-C), then rewrite
# the new, smaller table (rows 2-6).
$ws.Range("A2:C10").Clear()

# Row 1 (header) stays the same - no changes needed.

# Row 2
$ws.Cells.Item(2, 1).Value = " iAU_TC_ID_205"
$ws.Cells.Item(2, 2).Value = "@RegressionA Pre-Request Verify Elumina Login and Create Exam"
$ws.Cells.Item(2, 3).Value = "passed"

# Row 3
$ws.Cells.Item(3, 1).Value = " iAU_TC_ID_205"
$ws.Cells.Item(3, 2).Value = '@RegressionA Pre-Request "Validation of Delivery --> Add New Users"'
$ws.Cells.Item(3, 3).Value = "passed"

# Row 4
$ws.Cells.Item(4, 1).Value = "iAU_TC_ID_180"
$ws.Cells.Item(4, 2).Value = "@RegressionA Validation of Manage Delivery --> Edit user "
$ws.Cells.Item(4, 3).Value = "passed"

# Row 5
$ws.Cells.Item(5, 1).Value = "iAU_TC_ID_205"
$ws.Cells.Item(5, 2).Value = "@RegressionA Validation of Delivery --> Venue Summary "
$ws.Cells.Item(5, 3).Value = "passed"

# Row 6 (no value in column B)
$ws.Cells.Item(6, 1).Value = "iAU_TC_ID_206.,iAU_TC_ID_210.,iAU_TC_ID_209.,iAU_TC_ID_211 @RegressionA Validation of Delivery --> Live Dashboard "
$ws.Cells.Item(6, 3).Value = "timedOut"

$wb.Save()
